# Auto-generated edit script applying numeric corrections to Leve profit
# calculation columns (H-N) across multiple sheets, per the scheduled
# runner's recalculated market-price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 664.6667
$ws.Range("I4").Value = 664.6667
$ws.Range("K4").Value = 664.6667
$ws.Range("M4").Value = -550.6667
$ws.Range("H9").Value = 195.66667
$ws.Range("I9").Value = 102.166664
$ws.Range("J9").Value = 382.66666
$ws.Range("K9").Value = 102.166664
$ws.Range("L9").Value = 382.66666
$ws.Range("M9").Value = 66.833336
$ws.Range("N9").Value = -720.66666
$ws.Range("H96").Value = 1775
$ws.Range("I96").Value = 1045
$ws.Range("K96").Value = 3135
$ws.Range("M96").Value = -1762
$ws.Range("H97").Value = 4971
$ws.Range("J97").Value = 4587.778
$ws.Range("L97").Value = 13763.334
$ws.Range("N97").Value = -14755.334
$ws.Range("H99").Value = 184.16667
$ws.Range("I99").Value = 184.16667
$ws.Range("K99").Value = 552.50001
$ws.Range("M99").Value = 945.49999
$ws.Range("H113").Value = 4955.3
$ws.Range("I113").Value = 4510
$ws.Range("J113").Value = 5623.25
$ws.Range("K113").Value = 4510
$ws.Range("L113").Value = 5623.25
$ws.Range("M113").Value = -1256
$ws.Range("N113").Value = -12131.25
$ws.Range("H138").Value = 4918.106
$ws.Range("J138").Value = 4939.6943
$ws.Range("L138").Value = 14819.0829
$ws.Range("N138").Value = -25099.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2665.1667
$ws.Range("I88").Value = 1991
$ws.Range("K88").Value = 1991
$ws.Range("M88").Value = -1585
$ws.Range("H91").Value = 2665.1667
$ws.Range("I91").Value = 1991
$ws.Range("K91").Value = 1991
$ws.Range("M91").Value = -587
$ws.Range("H97").Value = 502.55554
$ws.Range("I97").Value = 516.6923
$ws.Range("J97").Value = 465.8
$ws.Range("K97").Value = 516.6923
$ws.Range("L97").Value = 465.8
$ws.Range("M97").Value = -20.69230000000005
$ws.Range("N97").Value = -1457.8
$ws.Range("H122").Value = 288686.9
$ws.Range("I122").Value = 387040.53
$ws.Range("J122").Value = 4554.222
$ws.Range("K122").Value = 1161121.59
$ws.Range("L122").Value = 13662.666
$ws.Range("M122").Value = -1158671.59
$ws.Range("N122").Value = -18562.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2197.6667
$ws.Range("J99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("N99").Value = -5496
$ws.Range("H134").Value = 3310.7273
$ws.Range("I134").Value = 2546.889
$ws.Range("K134").Value = 7640.667
$ws.Range("M134").Value = -5105.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 346.6
$ws.Range("I22").Value = 377.66666
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 377.66666
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -27.66665999999998
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 5114.037
$ws.Range("I31").Value = 5466.5386
$ws.Range("J31").Value = 4786.7144
$ws.Range("K31").Value = 5466.5386
$ws.Range("L31").Value = 4786.7144
$ws.Range("M31").Value = -5171.5386
$ws.Range("N31").Value = -5376.7144
$ws.Range("H34").Value = 5114.037
$ws.Range("I34").Value = 5466.5386
$ws.Range("J34").Value = 4786.7144
$ws.Range("K34").Value = 5466.5386
$ws.Range("L34").Value = 4786.7144
$ws.Range("M34").Value = -5264.5386
$ws.Range("N34").Value = -5190.7144
$ws.Range("H58").Value = 3821.76
$ws.Range("I58").Value = 2023.1177
$ws.Range("K58").Value = 2023.1177
$ws.Range("M58").Value = -1820.1177
$ws.Range("H86").Value = 12838.333
$ws.Range("I86").Value = 4762.25
$ws.Range("K86").Value = 4762.25
$ws.Range("M86").Value = -3639.25
$ws.Range("H89").Value = 12838.333
$ws.Range("I89").Value = 4762.25
$ws.Range("K89").Value = 23811.25
$ws.Range("M89").Value = -18195.25
$ws.Range("H102").Value = 43998
$ws.Range("J102").Value = 43998
$ws.Range("L102").Value = 43998
$ws.Range("N102").Value = -48866
$ws.Range("H109").Value = 37192.934
$ws.Range("J109").Value = 37192.934
$ws.Range("L109").Value = 37192.934
$ws.Range("N109").Value = -39272.934
$ws.Range("H122").Value = 884.7
$ws.Range("I122").Value = 878.2857
$ws.Range("J122").Value = 899.6667
$ws.Range("K122").Value = 2634.8571
$ws.Range("L122").Value = 2699.0001
$ws.Range("M122").Value = -184.8571000000002
$ws.Range("N122").Value = -7599.0001
$ws.Range("H132").Value = 2289.16
$ws.Range("J132").Value = 2749.5
$ws.Range("L132").Value = 8248.5
$ws.Range("N132").Value = -13308.5
$ws.Range("H136").Value = 3821.76
$ws.Range("I136").Value = 2023.1177
$ws.Range("K136").Value = 6069.3531
$ws.Range("M136").Value = -3519.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200.16667
$ws.Range("I2").Value = 37.8
$ws.Range("J2").Value = 316.14285
$ws.Range("K2").Value = 226.8
$ws.Range("L2").Value = 1896.8571
$ws.Range("M2").Value = -113.8
$ws.Range("N2").Value = -2122.8571
$ws.Range("H38").Value = 171.27272
$ws.Range("I38").Value = 141
$ws.Range("K38").Value = 423
$ws.Range("M38").Value = -76
$ws.Range("H117").Value = 1434
$ws.Range("J117").Value = 2297.5
$ws.Range("L117").Value = 6892.5
$ws.Range("N117").Value = -13776.5
$ws.Range("H132").Value = 2047
$ws.Range("J132").Value = 2047
$ws.Range("L132").Value = 18423
$ws.Range("N132").Value = -23483
$ws.Range("H137").Value = 6666.6
$ws.Range("I137").Value = 10000
$ws.Range("J137").Value = 5833.25
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 17499.75
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = -27699.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 18187.5
$ws.Range("J48").Value = 21375
$ws.Range("L48").Value = 21375
$ws.Range("N48").Value = -22345
$ws.Range("H74").Value = 15000
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -16872
$ws.Range("H77").Value = 15000
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -54360
$ws.Range("H97").Value = 1435.8966
$ws.Range("I97").Value = 1548.5238
$ws.Range("J97").Value = 1140.25
$ws.Range("K97").Value = 1548.5238
$ws.Range("L97").Value = 1140.25
$ws.Range("M97").Value = -1052.5238
$ws.Range("N97").Value = -2132.25
$ws.Range("H122").Value = 462176.5
$ws.Range("I122").Value = 74774.79
$ws.Range("J122").Value = 1004538.9
$ws.Range("K122").Value = 224324.37
$ws.Range("L122").Value = 3013616.7
$ws.Range("M122").Value = -221874.37
$ws.Range("N122").Value = -3018516.7
$ws.Range("H126").Value = 4994.7
$ws.Range("I126").Value = 4973.5
$ws.Range("K126").Value = 14920.5
$ws.Range("M126").Value = -12450.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3821.2144
$ws.Range("J46").Value = 4612.125
$ws.Range("L46").Value = 4612.125
$ws.Range("N46").Value = -4988.125
$ws.Range("H68").Value = 5000
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 5000
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H93").Value = 913.9167
$ws.Range("I93").Value = 440.77777
$ws.Range("K93").Value = 440.77777
$ws.Range("M93").Value = 807.2222300000001
$ws.Range("H100").Value = 2157.6
$ws.Range("I100").Value = 2197
$ws.Range("K100").Value = 2197
$ws.Range("M100").Value = -1656
$ws.Range("H127").Value = 93125
$ws.Range("J127").Value = 93125
$ws.Range("L127").Value = 93125
$ws.Range("N127").Value = -103045
$ws.Range("H132").Value = 5517.8667
$ws.Range("J132").Value = 6318.091
$ws.Range("L132").Value = 18954.273
$ws.Range("N132").Value = -24014.273
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("H141").Value = 70715
$ws.Range("J141").Value = 70715
$ws.Range("L141").Value = 70715
$ws.Range("N141").Value = -81075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3347.111
$ws.Range("I81").Value = 3202.75
$ws.Range("K81").Value = 6405.5
$ws.Range("M81").Value = -5344.5
$ws.Range("H84").Value = 3347.111
$ws.Range("I84").Value = 3202.75
$ws.Range("K84").Value = 32027.5
$ws.Range("M84").Value = -26723.5
$ws.Range("H132").Value = 3331.7778
$ws.Range("I132").Value = 2497.5
$ws.Range("K132").Value = 7492.5
$ws.Range("M132").Value = -4962.5
